$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Finished Projects")

# --- Remember the exact formatting of the two hyperlink cells (D5, D8) before
# --- touching anything, by stashing copies of their formats in a scratch area
# --- far away from the table (row 20) so it is unaffected by the row insert
# --- that happens at row 2 (Excel/this engine extends row-1 formatting down
# --- into a freshly inserted row 2, which would otherwise contaminate a
# --- nearby scratch cell on row 1).
$ws.Range("D8").Copy()
$ws.Range("F20").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("D5").Copy()
$ws.Range("G20").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# --- Insert a new row at row 2 (pushes existing rows 2-8 down to rows 3-9) ---
$ws.Rows.Item(2).Insert()

# Fill in the new row's data: No.=8, Project name, My role=Developer
$ws.Cells.Item(2, 1).Value = 8
$ws.Cells.Item(2, 2).Value = "License Plate Recognition using YOLOv10"
$ws.Cells.Item(2, 3).Value = "Developer"

# Copy formatting from the row below (old row 2, now shifted to row 3) so the
# new row matches the styling used throughout the rest of the table
$ws.Range("A3:D3").Copy()
$ws.Range("A2:D2").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# The role/notes column (C) on the new row uses a slightly different number
# format than the row that was copied from (General instead of #,##0)
$ws.Range("C2").NumberFormat = "#,##0"

# Restore the row height (PasteSpecial formats does not carry row height)
$ws.Rows.Item(2).RowHeight = 21.75

# --- Fix up the two hyperlinks, which do not automatically follow the cells
# --- they were attached to when the row was inserted above them. They need
# --- to move from D8 -> D9 and D5 -> D6.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D9"), "https://github.com/nntrivi2001/Face-recognition-with-GUI", "", "", "Coursework - Github")
$ws.Range("D9").Value = "Github"
$ws.Range("F21").Copy()
$ws.Range("D9").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Hyperlinks.Add($ws.Range("D6"), "https://github.com/nntrivi2001/SecureMemo", "", "", "Coursework - Github")
$ws.Range("D6").Value = "Github"
$ws.Range("G21").Copy()
$ws.Range("D6").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Clean up the scratch cells (the row-2 insert shifted our scratch row from 20 to 21)
$ws.Range("F21:G21").Clear()

# Match the final selected cell recorded in the workbook
$ws.Range("C3").Select()
